$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Seed" column (C) to make room for
# "Optimizer". Insert() shifts C:F -> D:G and keeps formulas/merges in sync.
$ws.Columns("C:C").Insert() | Out-Null

# Make the new column's width match its left neighbour ("batch size"),
# mirroring the target col min="2" max="3" grouping.
$ws.Columns("C:C").ColumnWidth = $ws.Columns("B:B").ColumnWidth

# New header
$ws.Range("C2").Value = "Optimizer"

# Existing rows (3-7) all used the Adam optimizer
$ws.Range("C3").Value = "Adam"
$ws.Range("C4").Value = "Adam"
$ws.Range("C5").Value = "Adam"
$ws.Range("C6").Value = "Adam"
$ws.Range("C7").Value = "Adam"

# New rows 8 and 9 - reuse row 7's formatting (batch size=1024 banded style)
$ws.Range("A7:G7").Copy() | Out-Null
$ws.Range("A8:A9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 8: batch size 1024, SGD optimizer, seed 0
$ws.Range("A8").Value = "CodeGPTPy"
$ws.Range("B8").Value = 1024
$ws.Range("C8").Value = "SGD"
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.673
$ws.Range("F8").Value = 0.689
$ws.Range("G8").Formula = "=E8-F8"

# Row 9: batch size 12244, Adam optimizer, seed 0
$ws.Range("A9").Value = "CodeGPTPy"
$ws.Range("B9").Value = 12244
$ws.Range("C9").Value = "Adam"
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.58
$ws.Range("F9").Value = 0.675
$ws.Range("G9").Formula = "=E9-F9"

# Match the saved selection from the authored edit
$ws.Range("F10").Select() | Out-Null
